$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# G2 was empty - fill in the sound tag for row 2 (matches A2 file name)
$ws.Range("G2").Value = "[sound:0003_ഒപ്പം_01.mp3]"

# H2:H11 were empty - fill in the image tag for each row using the
# same file-name identifier stored in column A of that row.
for ($row = 2; $row -le 11; $row++) {
    $fileName = $ws.Cells.Item($row, 1).Value2
    $ws.Cells.Item($row, 8).Value = '<img src="' + $fileName + '.jpg">'
}
